$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing, so numeric-looking strings
# (e.g. "1.00", "636.68") are preserved as text instead of being parsed as numbers.
# The format is reset back to Normal afterwards so no stray style is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '79.696.62'
$ws.Range("E2").Value = '  +4.25%  '
$ws.Range("D3").Value = '3.214.29'
$ws.Range("E3").Value = '  +5.52%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '205.75'
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("D6").Value = '636.68'
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '0.243'
$ws.Range("E8").Value = '  +17.65%  '
$ws.Range("E9").Value = '  +6.05%  '
$ws.Range("D10").Value = '3.218.48'
$ws.Range("E10").Value = '  +5.59%  '
$ws.Range("D11").Value = '0.590'
$ws.Range("E11").Value = '  +34.35%  '
$ws.Range("D12").Value = '0.165'
$ws.Range("E12").Value = '  +2.99%  '
$ws.Range("D13").Value = '5.49'
$ws.Range("E13").Value = '  +4.22%  '
$ws.Range("D14").Value = '0.0000233'
$ws.Range("E14").Value = '  +20.63%  '
$ws.Range("D15").Value = '3.812.72'
$ws.Range("E15").Value = '  +5.61%  '
$ws.Range("D16").Value = '31.84'
$ws.Range("E16").Value = '  +8.52%  '
$ws.Range("D17").Value = '79.684.75'
$ws.Range("E17").Value = '  +4.29%  '
$ws.Range("D18").Value = '3.219.64'
$ws.Range("E18").Value = '  +5.56%  '
$ws.Range("D19").Value = '14.50'
$ws.Range("E19").Value = '  +6.88%  '
$ws.Range("D20").Value = '3.02'
$ws.Range("E20").Value = '  +29.69%  '
$ws.Range("D21").Value = '9.30'
$ws.Range("E21").Value = '  +2.66%  '
$ws.Range("D22").Value = '430.20'
$ws.Range("E22").Value = '  +14.25%  '
$ws.Range("D23").Value = '5.12'
$ws.Range("E23").Value = '  +17.09%  '
$ws.Range("D24").Value = '3.384.95'
$ws.Range("E24").Value = '  +5.61%  '
$ws.Range("D25").Value = '11.19'
$ws.Range("E25").Value = '  +13.08%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '77.03'
$ws.Range("E26").Value = '  +4.60%  '
$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D27").Value = '4.71'
$ws.Range("E27").Value = '  +6.72%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("E29").Value = '  +7.39%  '
$ws.Range("D30").Value = '9.04'
$ws.Range("E30").Value = '  +8.52%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("D32").Value = '1.48'
$ws.Range("E32").Value = '  +4.42%  '
$ws.Range("D33").Value = '527.76'
$ws.Range("E33").Value = '  +4.97%  '
$ws.Range("D34").Value = '2.00'
$ws.Range("E34").Value = '  +2.00%  '
$ws.Range("E35").Value = '  +21.14%  '
$ws.Range("D36").Value = '23.11'
$ws.Range("E36").Value = '  +10.82%  '
$ws.Range("D37").Value = '0.122'
$ws.Range("E37").Value = '  +14.90%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("D39").Value = '0.409'
$ws.Range("E39").Value = '  +5.68%  '
$ws.Range("D40").Value = '165.14'
$ws.Range("E40").Value = '  +1.24%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '192.65'
$ws.Range("E42").Value = '  +0.64%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '5.52'
$ws.Range("E44").Value = '  +6.59%  '
$ws.Range("E45").Value = '  +2.84%  '
$ws.Range("D46").Value = '1.79'
$ws.Range("E46").Value = '  +7.53%  '
$ws.Range("E47").Value = '  +3.41%  '
$ws.Range("D48").Value = '43.62'
$ws.Range("E48").Value = '  +3.51%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '25.87'
$ws.Range("E49").Value = '  +14.88%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = '2.55'
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").Value = '0.637'
$ws.Range("E51").Value = '  +4.27%  '

# Reset column D style back to Normal (removes the temporary text-format style index)
$dRange.Style = "Normal"

